$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" column (O) mirroring the existing "2020" column (N):
# copy each N-row's formatting into the matching O-row, then set the new values.

$ws.Range("N3").Copy()
$ws.Range("O3").PasteSpecial(-4122)

$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("O4").Value = 2021

$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)

$ws.Range("N6").Copy()
$ws.Range("O6").PasteSpecial(-4122)
$ws.Range("O6").Value = 1860

$ws.Range("N7").Copy()
$ws.Range("O7").PasteSpecial(-4122)
$ws.Range("O7").Value = 1

$ws.Range("N8").Copy()
$ws.Range("O8").PasteSpecial(-4122)
$ws.Range("O8").Value = 510

$ws.Range("N9").Copy()
$ws.Range("O9").PasteSpecial(-4122)
$ws.Range("O9").Value = 178

$ws.Range("N10").Copy()
$ws.Range("O10").PasteSpecial(-4122)
$ws.Range("O10").Value = 821

$excel.CutCopyMode = $false

# Match the author's final selection state.
$ws.Range("P9").Select()
